$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trading-day dates for rows 182..190 (row 181 already holds "20-09-2021").
$dates = @(
    "21-09-2021",
    "22-09-2021",
    "23-09-2021",
    "24-09-2021",
    "27-09-2021",
    "28-09-2021",
    "29-09-2021",
    "30-09-2021",
    "01-10-2021"
)

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = 182 + $i
    $text = $dates[$i]
    if ($row -eq 190) {
        # "01-10-2021" is ambiguous (day <= 12) and Excel's smart text
        # parsing would otherwise convert it into a date serial number.
        # Route it through a text formula and paste-special as values so
        # the cell lands as a plain shared-string, with no left-over
        # number-format/style applied to the cell.
        $ws.Cells.Item($row, 1).Formula = '="' + $text + '"'
        $ws.Cells.Item($row, 1).Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4163)
    }
    else {
        $ws.Cells.Item($row, 1).Value = $text
    }
}

# Row-by-row numeric data (columns B..G) for rows 181..190.
$data = @(
    @{ Row = 181; B = 80000;  C = 1.5; D = 1.5; E = 1.5; F = 5; G = 1.5 },
    @{ Row = 182; B = 110000; C = 1.5; D = 1.5; E = 1.5; F = 4; G = 1.5 },
    @{ Row = 183; B = 40000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ Row = 184; B = 180000; C = 1.5; D = 1.5; E = 1.5; F = 4; G = 1.5 },
    @{ Row = 185; B = 70000;  C = 1.5; D = 1.5; E = 1.5; F = 4; G = 1.5 },
    @{ Row = 186; B = 85000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ Row = 187; B = 85000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ Row = 188; B = 80000;  C = 1.5; D = 1.5; E = 1.5; F = 3; G = 1.5 },
    @{ Row = 189; B = 20000;  C = $null; D = $null; E = $null; F = 2; G = 1.5 },
    @{ Row = 190; B = $null; C = $null; D = $null; E = $null; F = $null; G = 1.5 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    if ($null -ne $entry.B) { $ws.Cells.Item($r, 2).Value = $entry.B }
    if ($null -ne $entry.C) { $ws.Cells.Item($r, 3).Value = $entry.C }
    if ($null -ne $entry.D) { $ws.Cells.Item($r, 4).Value = $entry.D }
    if ($null -ne $entry.E) { $ws.Cells.Item($r, 5).Value = $entry.E }
    if ($null -ne $entry.F) { $ws.Cells.Item($r, 6).Value = $entry.F }
    if ($null -ne $entry.G) { $ws.Cells.Item($r, 7).Value = $entry.G }
}
